# Insert a new weekly record as row 558, pushing the existing rows 558:580
# down to 559:581 (dimension grows from R580 to R581).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("558:558").Insert()

$row = 558

$ws.Cells.Item($row, 1).Value = 10
$ws.Cells.Item($row, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($row, 3).Value = "La Araucanía"
$ws.Cells.Item($row, 4).Value = 45008
$ws.Cells.Item($row, 5).Value = 9
$ws.Cells.Item($row, 6).Value = 100112037
$ws.Cells.Item($row, 7).Value = "Cebollín"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 110
$ws.Cells.Item($row, 11).Value = 7000
$ws.Cells.Item($row, 12).Value = 7000
$ws.Cells.Item($row, 13).Value = 7000
$ws.Cells.Item($row, 14).Value = "$/docena de paquetes"
$ws.Cells.Item($row, 15).Value = "Provincia de Cautín"
$ws.Cells.Item($row, 16).Value = 583
$ws.Cells.Item($row, 17).Value = 12
$ws.Cells.Item($row, 18).Value = "Hortaliza"
